$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "This is an essential problem to the CG bias. We admit that the current simulations cannot represent the HST data. The main reason of failure fitting two Sersic is that some galaxies cannot model by that. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This is indeed a relevant limitations of the CG bias analysis. We admit that the current simulations cannot fully represent the HST data. The main reason of failure fitting two Sersic models is that some galaxies cannot be modelled by that. ",
    2)

$d.Content.Find.Execute(
    "It requires some further collaboration in future with other groups,  including the image simulation and shear measurements. For the current problem, we perform two additional tests: 1) in the simulations, we also fit the galaxies with one Sersic. The estimation has some extra error, but can reflect the main properties of the bias. 2) We also fit the HST data with two Sersics, in the successful fitting (about 70%), the estimation shows similar relation with the colour of the galaxy, although the scatter is larger.  We add a small part to explain this point at the end of section 4.1.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "It requires some further collaboration in the future with other groups,  including the image simulation and shear measurements. For the current problem, we perform two additional tests: 1) in the simulations, we also fit the galaxies with one Sersic component. The estimation has some extra scatter, but can reflect the main properties of the bias. 2) We also fit the HST data with two Sersic components, in the successful fitting (about 70%), the estimation shows a similar relation with the colour of the galaxy, although the scatter is larger.  We add a small part to explain this point at the end of section 4.1.",
    2)
